# Generate Report for Handback
#
# Populates the "Latest Target File" (F) and "Latest Handback File" (G)
# columns -- previously empty -- for both language sheets (zh-cn, de-de),
# flips the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US", and records the handback timestamps in
# the "Latest Handback DateTime" column (H).

$wb = $excel.ActiveWorkbook

function Add-HandbackLink($ws, $cellRef, $address, $displayText) {
    # Hyperlinks.Add both sets the cell text (TextToDisplay) and wires up the
    # external relationship; re-applying the built-in "HyperLink" style keeps
    # the same look'n'feel as the pre-existing Handoff-file hyperlinks.
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText) | Out-Null
    $ws.Range($cellRef).Style = "HyperLink"
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 -> b85dda28-...  Row 3 -> edc33579-...
Add-HandbackLink $wsZh "F2" "https://github.com/OpenLocalizationTest/oltest/blob/93cf9d62869571819dc7492d3511edc70957fcac/e2e/b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.md" "b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.md"
Add-HandbackLink $wsZh "G2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/93cf9d62869571819dc7492d3511edc70957fcac/ol-handback/OpenLocalizationTest/oltest/xinjiang/b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.3bfa9c417027403a4f17b00a95c1e886900af3fe.zh-cn.xlf" "b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.3bfa9c417027403a4f17b00a95c1e886900af3fe.zh-cn.xlf"

Add-HandbackLink $wsZh "F3" "https://github.com/OpenLocalizationTest/oltest/blob/93cf9d62869571819dc7492d3511edc70957fcac/e2e/edc33579-ba42-432e-a4d2-a999004f8f1d.md" "edc33579-ba42-432e-a4d2-a999004f8f1d.md"
Add-HandbackLink $wsZh "G3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/93cf9d62869571819dc7492d3511edc70957fcac/ol-handback/OpenLocalizationTest/oltest/xinjiang/edc33579-ba42-432e-a4d2-a999004f8f1d.ebd44a9eda907cf83596255cbd5686dfffe7f378.zh-cn.xlf" "edc33579-ba42-432e-a4d2-a999004f8f1d.ebd44a9eda907cf83596255cbd5686dfffe7f378.zh-cn.xlf"

# Status -> handed back, in sync with en-US
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime
$wsZh.Range("H2").Value = "2016-03-24 03:04:05"
$wsZh.Range("H3").Value = "2016-03-24 03:04:05"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Add-HandbackLink $wsDe "F2" "https://github.com/OpenLocalizationTest/oltest/blob/93cf9d62869571819dc7492d3511edc70957fcac/e2e/b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.md" "b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.md"
Add-HandbackLink $wsDe "G2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ebab86994c1e47dbaa23ca132f44b8ceeda73f0/ol-handback/OpenLocalizationTest/oltest/xinjiang/b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.3bfa9c417027403a4f17b00a95c1e886900af3fe.de-de.xlf" "b85dda28-b103-4b8b-87aa-f2c9a1fc5fe9.3bfa9c417027403a4f17b00a95c1e886900af3fe.de-de.xlf"

Add-HandbackLink $wsDe "F3" "https://github.com/OpenLocalizationTest/oltest/blob/93cf9d62869571819dc7492d3511edc70957fcac/e2e/edc33579-ba42-432e-a4d2-a999004f8f1d.md" "edc33579-ba42-432e-a4d2-a999004f8f1d.md"
Add-HandbackLink $wsDe "G3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ebab86994c1e47dbaa23ca132f44b8ceeda73f0/ol-handback/OpenLocalizationTest/oltest/xinjiang/edc33579-ba42-432e-a4d2-a999004f8f1d.ebd44a9eda907cf83596255cbd5686dfffe7f378.de-de.xlf" "edc33579-ba42-432e-a4d2-a999004f8f1d.ebd44a9eda907cf83596255cbd5686dfffe7f378.de-de.xlf"

# Status -> handed back, in sync with en-US
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime (de-de finished a little later than zh-cn)
$wsDe.Range("H2").Value = "2016-03-24 03:04:19"
$wsDe.Range("H3").Value = "2016-03-24 03:04:19"

Write-Output "Handback report generated."
